$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Reorder slides: the "Q&A" slide (currently slide 19) moves up to become
#    slide 17, pushing "Plans For Next Sprint" and the chess-pools slide down
#    by one position each.
# ---------------------------------------------------------------------------
$p.Slides.Item(19).MoveTo(17)

# ---------------------------------------------------------------------------
# 2) Rework the agenda bullets on slide 3 ("Overview of Presentation (2/2)")
#    so the order becomes:
#       Demonstration of Application
#       Q&A
#       Plans For Next Sprint
#    and the last line is split back into two runs ("Plans " / "For Next Sprint").
# ---------------------------------------------------------------------------
$agendaSlide = $p.Slides.Item(3)
$contentShape = $agendaSlide.Shapes.Item(2)
$tr = $contentShape.TextFrame.TextRange
$tr.Text = "Demonstration of Application`rQ&A`rPlans For Next Sprint"

$lastPara = $tr.Paragraphs(3, 1)
$firstRun = $tr.Characters($lastPara.Start, 6)
$firstRun.Text = "Plans "
